$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style (bold/centered/bordered) as the other header cells
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Data cells F2 / F3 - plain values (no special style, matching E2/E3)
$ws.Range("F2").Value = "2021-10-05 13:39:38.320660"
$ws.Range("F3").Value = "2021-10-05 13:39:38.320673"
